$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "BiBBS_CohortInfo.pregnancy.recruitment_date"
$ws.Range("A19").Select()
